# Update countries & provincias Spain
# Applies the 30-Jun-2020 13:32 data refresh on top of the 12:15 snapshot:
#  - Country rows whose rank (sorted by "Casos totales" desc) moved up by one
#    swap their country-name label with the row above (Madagascar/Paraguay,
#    Malaui/Hong Kong, Dominica/Fiyi, Groenlandia/Islas Malvinas).
#  - Updated case counters (Casos totales/Nuevos casos/Casos activos/
#    Recuperados/Casos criticos/Muertes hoy/Muertes) for the affected rows.
#  - The "Datos actualizados" footer timestamp moves from 12:15 to 13:32.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country label swaps (column A) caused by the re-sort ---------------
$ws.Range("A107").Value = "Madagascar"
$ws.Range("A108").Value = "Paraguay"
$ws.Range("A109").Value = "Mali"
$ws.Range("A110").Value = "Nicaragua"

$ws.Range("A124").Value = "Malaui"
$ws.Range("A125").Value = "Hong Kong"
$ws.Range("A126").Value = "Benin"
$ws.Range("A127").Value = "Tunez"
$ws.Range("A128").Value = "Cabo Verde"

$ws.Range("A205").Value = "Dominica"
$ws.Range("A206").Value = "Fiyi"

$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"

# --- Updated statistics (columns B-H) ------------------------------------
$ws.Range("B7").Value = 568473
$ws.Range("C7").Value = 937
$ws.Range("E7").Value = 215898

$ws.Range("B13").Value = 227662
$ws.Range("C13").Value = 2457
$ws.Range("D13").Value = 188758
$ws.Range("E13").Value = 28087
$ws.Range("G13").Value = 147
$ws.Range("H13").Value = 10817

$ws.Range("B46").Value = 31714
$ws.Range("C46").Value = 62
$ws.Range("E46").Value = 652

$ws.Range("B48").Value = 26970
$ws.Range("C48").Value = 388
$ws.Range("D48").Value = 19050
$ws.Range("E48").Value = 6269
$ws.Range("G48").Value = 17
$ws.Range("H48").Value = 1651

$ws.Range("B64").Value = 13564
$ws.Range("C64").Value = 316
$ws.Range("D64").Value = 3194
$ws.Range("E64").Value = 10341

$ws.Range("D74").Value = 5569
$ws.Range("E74").Value = 2704
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 25

$ws.Range("B78").Value = 6793
$ws.Range("C78").Value = 95
$ws.Range("D78").Value = 4431
$ws.Range("E78").Value = 2250
$ws.Range("G78").Value = 4
$ws.Range("H78").Value = 112

$ws.Range("B107").Value = 2214
$ws.Range("C107").Value = 76
$ws.Range("D107").Value = 994
$ws.Range("E107").Value = 1200
$ws.Range("H107").Value = 20

$ws.Range("B108").Value = 2191
$ws.Range("D108").Value = 1080
$ws.Range("E108").Value = 1095
$ws.Range("H108").Value = 16

$ws.Range("B109").Value = 2173
$ws.Range("D109").Value = 1447
$ws.Range("E109").Value = 611
$ws.Range("H109").Value = 115

$ws.Range("B110").Value = 2170
$ws.Range("D110").Value = 1238
$ws.Range("E110").Value = 858
$ws.Range("H110").Value = 74

$ws.Range("B124").Value = 1224
$ws.Range("C124").Value = 72
$ws.Range("D124").Value = 260
$ws.Range("E124").Value = 950
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 14

$ws.Range("B125").Value = 1204
$ws.Range("D125").Value = 1105
$ws.Range("E125").Value = 92
$ws.Range("H125").Value = 7

$ws.Range("B126").Value = 1199
$ws.Range("C126").Value = 12
$ws.Range("D126").Value = 333
$ws.Range("E126").Value = 845
$ws.Range("G126").Value = 2
$ws.Range("H126").Value = 21

$ws.Range("B127").Value = 1172
$ws.Range("D127").Value = 1029
$ws.Range("E127").Value = 93
$ws.Range("H127").Value = 50

$ws.Range("B128").Value = 1165
$ws.Range("D128").Value = 608
$ws.Range("E128").Value = 545
$ws.Range("H128").Value = 12

$ws.Range("B136").Value = 962
$ws.Range("C136").Value = 3
$ws.Range("D136").Value = 838
$ws.Range("E136").Value = 71

$ws.Range("D150").Value = 640
$ws.Range("E150").Value = 21

# --- Footer timestamp -----------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Junio de 2020 a las 13:32"
